# Added Map to reports (City tab)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: only the username changes
$ws.Range("A7").Value = "safafgas"

# Row 8
$ws.Range("A8").Value = "zMaxShap"
$ws.Range("B8").Value = "asfashjk"
$ws.Range("C8").Value = "Sman95@gmail.com"

# Row 9
$ws.Range("A9").Value = "asfda"
$ws.Range("B9").Value = "fasfas"
$ws.Range("C9").Value = "asffa"

# Row 10
$ws.Range("A10").Value = "fasfa"
$ws.Range("B10").Value = "asfafa"
$ws.Range("C10").Value = "fasf"

# Row 11
$ws.Range("A11").Value = "gsdsdg"
$ws.Range("B11").Value = "sgdsg"
$ws.Range("C11").Value = "sdgsg"

# Row 12
$ws.Range("A12").Value = "dfhdgf"
$ws.Range("B12").Value = "gfawgfafs"
$ws.Range("C12").Value = "asfasfas"

# Row 13
$ws.Range("A13").Value = "sfafaf"
$ws.Range("B13").Value = "asfafs"
$ws.Range("C13").Value = "Sman95@gmail.com"

# Row 14
$ws.Range("A14").Value = "asfaf"
$ws.Range("B14").Value = "asasf"
$ws.Range("C14").Value = "afs"

# Row 15
$ws.Range("A15").Value = "MaxShapira"
$ws.Range("B15").Value = "asfnajskflh"
$ws.Range("C15").Value = "Sman95@gmail.com"

# Row 16
$ws.Range("A16").Value = "zMaxShap"
$ws.Range("B16").Value = "fjashdfkas"
$ws.Range("C16").Value = "afafgasf"

# Row 17
$ws.Range("A17").Value = "zMaxShap"
$ws.Range("B17").Value = "hasikfl"
$ws.Range("C17").Value = "sdgfjksglfs"

# Row 18
$ws.Range("A18").Value = "asfjakf"
$ws.Range("B18").Value = "dsafjaskfl"
$ws.Range("C18").Value = "Sman95@gmail.com"

# Row 19
$ws.Range("A19").Value = "asfafsjkl"
$ws.Range("B19").Value = "aqjisfaj"
$ws.Range("C19").Value = "Sman95@gmail.com"

# Row 20 (C20 stays empty)
$ws.Range("A20").Value = "xcvxv"
$ws.Range("B20").Value = "cxvxv"

# Row 21
$ws.Range("A21").Value = "fasfas"
$ws.Range("B21").Value = "fsafasf"
$ws.Range("C21").Value = "asfa"

# Row 22
$ws.Range("A22").Value = "fsafaf"
$ws.Range("B22").Value = "sfasfa"
$ws.Range("C22").Value = "sfafs"

# Update the sheet's selection to match the committed view
$ws.Range("A9:C9").Select() | Out-Null
